$wb = $excel.ActiveWorkbook
$wsCompare = $wb.Worksheets.Item("compare")
$wsHistory = $wb.Worksheets.Item("Change History")

# Update header row: G1 gets the new "Portfoio" label (matches F1's source
# header semantics), and B1 becomes the combined "Numerai  Portfoio" label.
# (G1 is written first so the shared-string table picks up "Portfoio"
# before "Numerai  Portfoio", matching the workbook's string order.)
$wsCompare.Range("G1").Value = "Portfoio"
$wsCompare.Range("B1").Value = "Numerai  Portfoio"

# Replace the placeholder weight column (G2:G347) with the real Numerai
# portfolio weights looked up from the source data.
$gVals = @(0.050299999999999997,0.81299999999999994,0.2079,0.29459999999999997,1.0092000000000001,0.76180000000000003,0.19270000000000001,0.5716,0.1066,0,0.089200000000000002,0.21440000000000001,0,0.66300000000000003,0.70640000000000003,0.090999999999999998,0.4924,0.37609999999999999,0.47839999999999999,0.038600000000000002,0,0,0.20150000000000001,0.0877,0.42449999999999999,0.2379,0,0.58069999999999999,0,0.099900000000000003,0.2157,0.039399999999999998,0,0.1512,0,0,1.0468999999999999,1.1507000000000001,0.1477,0.14119999999999999,0.087300000000000003,1.9424999999999999,0.73060000000000003,0.056000000000000001,0.8518,0.44729999999999998,0.18479999999999999,0.41520000000000001,0.014,0.54479999999999995,0.43209999999999998,0.2303,0,0,0,0.1239,0.38940000000000002,0.16669999999999999,0.14499999999999999,0,0.1081,0.14460000000000001,0.27189999999999998,0.27639999999999998,0,0.59560000000000002,0.622,0.16059999999999999,0.49299999999999999,0.35039999999999999,0,0.049200000000000001,0,0,0.12770000000000001,0.34649999999999997,0.99450000000000005,0.55579999999999996,1.7201,0,0.95109999999999995,0,0,0.34060000000000001,0.13539999999999999,0,0,0.86529999999999996,0.12039999999999999,1.5428999999999999,0,0,1.7658,0.39529999999999998,0,0.128,0,1.4851000000000001,0.1552,0.088900000000000007,0,0.10299999999999999,0.41470000000000001,0.1237,0.61260000000000003,1.3434999999999999,0.97840000000000005,0.17699999999999999,0.21579999999999999,0.44919999999999999,0.31080000000000002,0,0,0.46879999999999999,0.47889999999999999,0.22869999999999999,0,0.93730000000000002,0.153,0,0.40160000000000001,0.078600000000000003,0,0,0.034599999999999999,0.17879999999999999,0.1517,0.13950000000000001,0.43990000000000001,0.21060000000000001,0.40620000000000001,0.088900000000000007,0.4229,0,0.1739,0.55530000000000002,0.13539999999999999,0.6613,0.44569999999999999,0.20039999999999999,0.16200000000000001,0,0.23630000000000001,0.36899999999999999,0,0.23150000000000001,0.12479999999999999,1.0329999999999999,0.38419999999999999,0,0.4098,0.15379999999999999,0.42049999999999998,0.096500000000000002,1.3638999999999999,0,0,0,0.55959999999999999,0,0.309,0.090999999999999998,0,0.13500000000000001,0.20569999999999999,0,0.081000000000000003,0.74180000000000001,0.091700000000000004,0.0872,0.65359999999999996,2.1173999999999999,1.8093999999999999,0,0,0.48080000000000001,1.1579999999999999,0.15909999999999999,0,0.79510000000000003,0,0,0,0.26279999999999998,1.2236,0.31340000000000001,0,0.40179999999999999,0.097900000000000001,0,0,0.12620000000000001,0.45650000000000002,0.1361,0,0.29770000000000002,0.1153,0,0.071599999999999997,0.40500000000000003,0,0.33500000000000002,0.36699999999999999,0.26490000000000002,0.2031,0.84399999999999997,0.67249999999999999,0.1452,0.1143,0.28289999999999998,0.42830000000000001,0,0.015599999999999999,0.11609999999999999,0.027199999999999998,0.23849999999999999,1.1336999999999999,0,2.5878000000000001,0.1439,0.252,0.69489999999999996,0.2611,0.28289999999999998,0,0.16120000000000001,0.3886,0,0.33779999999999999,0,0.2051,0.26860000000000001,0.041700000000000001,0,0,0.40189999999999998,0,0,0.2853,0.1308,0,0,0.19070000000000001,0.58069999999999999,0.1258,1.8292999999999999,0.1081,0.13589999999999999,0.0974,0.088700000000000001,0.060199999999999997,0,0.18720000000000001,1.1234999999999999,0.1419,0,0,0.093399999999999997,0.30120000000000002,0.65690000000000004,0,0,0.1857,0,0.1188,0,0.2407,0.18629999999999999,0.44829999999999998,0,0.0241,0.1148,0.21890000000000001,0.47810000000000002,0.1343,0.1769,0,0,0.1331,0.1951,0.091300000000000006,0.1391,0.1069,0,0.5151,0.28199999999999997,0.53939999999999999,0.8468,0,0.12,0.1288,0.1105,0,0.0361,0,0.090499999999999997,0.47820000000000001,0.070599999999999996,0,0.29039999999999999,0,0,0.2205,0,0.59560000000000002,1.1167,0.1226,0.52090000000000003,0,0,0,0.4032,1.8391999999999999,0,0.14050000000000001,0.54400000000000004,0,0.27439999999999998,0.73409999999999997,0.0436,0.11609999999999999,0,0.1011,0.48659999999999998,0,0,0,0,1.1866000000000001,0.78469999999999995,0,0.099500000000000005,0.10920000000000001,0.34889999999999999,0.1308,0.2364,0.092600000000000002,0.71699999999999997,0.3488,0,0.091499999999999998,0.1115,0.095299999999999996,0.34589999999999999,0.40579999999999999,0.4592)
for ($i = 0; $i -lt $gVals.Length; $i++) {
    $wsCompare.Cells.Item($i + 2, 7).Value = $gVals[$i]
}

# Widen column B now that it holds the longer combined header text.
$wsCompare.Columns.Item(2).ColumnWidth = 25.45

# The lookup table now spans the whole F:G data range, so repoint the
# hidden _FilterDatabase name and the sheet's AutoFilter to match.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=compare!`$F`$1:`$G`$347"
$wsCompare.AutoFilterMode = $false
$wsCompare.Range("F1:G347").AutoFilter()

# Change History becomes the active/visible tab instead of compare.
$wsHistory.Activate()
